$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (A empty, B/C held the "5817372 - Simone de Fátima Medeiros
# Sampaio" docente string under "Docentes responsáveis:") is removed entirely;
# Excel shifts every row below it up by one (also shifting row heights).
$ws.Rows(13).Delete()

# After the shift, a handful of cells need their text updated in place.

# Row 10 ("Objetivos:") body text is replaced by the docente string.
$ws.Range("B10").Value = "5817372 - Simone de Fátima Medeiros Sampaio"
$ws.Range("C10").Value = "5817372 - Simone de Fátima Medeiros Sampaio"

# Row 13 ("Programa resumido:") body text becomes "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 ("Programa:") body text becomes the activation date "01/01/2018".
# A plain .Value assignment of that string gets auto-converted to a date
# serial by Excel's type sniffing, so instead pull the literal text value
# from the existing "Ativacao:" cell (B8, already stored as text) and then
# restore the row's normal wrap-text format (copied from B14/C14) on top.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# Row 18 ("Método:") body text is replaced by the docente string.
$ws.Range("B18").Value = "5817372 - Simone de Fátima Medeiros Sampaio"
$ws.Range("C18").Value = "5817372 - Simone de Fátima Medeiros Sampaio"
